# Updating the models with the consideration of the PVPP as feature
#
# This script rewrites the clear/cloudy sky GHI/DNI/DHI figures on both the
# "Daily" and "Hourly" sheets to reflect the re-run model output.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Daily": single data row (row 2), columns G:L
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily")

$daily.Range("G2").Value = 3190.89
$daily.Range("H2").Value = 6465.37
$daily.Range("I2").Value = 765.75
$daily.Range("J2").Value = 3190.89
$daily.Range("K2").Value = 6145.72
$daily.Range("L2").Value = 768.08

# ---------------------------------------------------------------------------
# Sheet "Hourly": rows 9-19, columns H:M
# ---------------------------------------------------------------------------
$hourly = $wb.Worksheets.Item("Hourly")

# Row 9 (hour 7)
$hourly.Range("H9").Value = 8.210000000000001
$hourly.Range("I9").Value = 61.04
$hourly.Range("K9").Value = 8.210000000000001
$hourly.Range("L9").Value = 16.08
$hourly.Range("M9").Value = 7.26

# Row 10 (hour 8)
$hourly.Range("H10").Value = 115.07
$hourly.Range("I10").Value = 445.95
$hourly.Range("J10").Value = 50.46
$hourly.Range("K10").Value = 115.07
$hourly.Range("L10").Value = 409.2
$hourly.Range("M10").Value = 46.13

# Row 11 (hour 9)
$hourly.Range("H11").Value = 262.96
$hourly.Range("I11").Value = 647.1
$hourly.Range("J11").Value = 74.09
$hourly.Range("K11").Value = 262.96
$hourly.Range("L11").Value = 627.99
$hourly.Range("M11").Value = 71.81999999999999

# Row 12 (hour 10)
$hourly.Range("H12").Value = 390.13
$hourly.Range("I12").Value = 744.01
$hourly.Range("J12").Value = 87.87
$hourly.Range("K12").Value = 390.13
$hourly.Range("L12").Value = 727.55
$hourly.Range("M12").Value = 87.51000000000001

# Row 13 (hour 11)
$hourly.Range("H13").Value = 477.02
$hourly.Range("I13").Value = 792.88
$hourly.Range("J13").Value = 95.61
$hourly.Range("K13").Value = 477.02
$hourly.Range("L13").Value = 768.65
$hourly.Range("M13").Value = 100.47

# Row 14 (hour 12)
$hourly.Range("H14").Value = 513.34
$hourly.Range("I14").Value = 810.5700000000001
$hourly.Range("J14").Value = 98.55
$hourly.Range("K14").Value = 513.34
$hourly.Range("L14").Value = 781.96
$hourly.Range("M14").Value = 106.43

# Row 15 (hour 13)
$hourly.Range("H15").Value = 495.14
$hourly.Range("I15").Value = 801.95
$hourly.Range("J15").Value = 97.08
$hourly.Range("K15").Value = 495.14
$hourly.Range("L15").Value = 775.52
$hourly.Range("M15").Value = 103.43

# Row 16 (hour 14)
$hourly.Range("H16").Value = 424.37
$hourly.Range("I16").Value = 764.67
$hourly.Range("J16").Value = 91.02
$hourly.Range("K16").Value = 424.37
$hourly.Range("L16").Value = 745.62
$hourly.Range("M16").Value = 92.41

# Row 17 (hour 15)
$hourly.Range("H17").Value = 309.12
$hourly.Range("I17").Value = 687.13
$hourly.Range("J17").Value = 79.48999999999999
$hourly.Range("K17").Value = 309.12
$hourly.Range("L17").Value = 668.99
$hourly.Range("M17").Value = 77.95999999999999

# Row 18 (hour 16)
$hourly.Range("H18").Value = 165.06
$hourly.Range("I18").Value = 533.04
$hourly.Range("J18").Value = 59.96
$hourly.Range("K18").Value = 165.06
$hourly.Range("L18").Value = 509.2

# Row 19 (hour 17)
$hourly.Range("H19").Value = 30.47
$hourly.Range("I19").Value = 177.03
$hourly.Range("K19").Value = 30.47
$hourly.Range("L19").Value = 114.97
